$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same table of events and
# need their "想去人数" (F column) counts refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 589
    $ws.Range("F3").Value = 3679
    $ws.Range("F4").Value = 102
    $ws.Range("F5").Value = 707
}
